$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values after repulling / pushing all data
# and recomputing means. Row 13 and row 18 are unchanged.
$updates = @{
    2  = -4
    3  = -8
    4  = 1
    5  = -1
    6  = -1
    7  = -1
    8  = -4
    9  = -2
    10 = -4
    11 = 1
    12 = -7
    14 = -4
    15 = -8
    16 = -1
    17 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
